$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = "2024-07-22 21:13:19"
$ws.Cells.Item(4, 2).Value = 19
$ws.Cells.Item(4, 3).Value = 15
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 7
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0.001
$ws.Cells.Item(4, 10).Value = 0.05
$ws.Cells.Item(4, 11).Value = 0.003
$ws.Cells.Item(4, 12).Value = 100
$ws.Cells.Item(4, 13).Value = 500
$ws.Cells.Item(4, 14).Value = 10
$ws.Cells.Item(4, 15).Value = 6
$ws.Cells.Item(4, 16).Value = 3
$ws.Cells.Item(4, 17).Value = 1000
$ws.Cells.Item(4, 18).Value = 5
$ws.Cells.Item(4, 19).Value = 1
$ws.Cells.Item(4, 20).Value = 20
$ws.Cells.Item(4, 21).Value = 0.7894736842105263
$ws.Cells.Item(4, 22).Value = "./Data/Electromecanica.xlsx"
$ws.Cells.Item(4, 23).Value = 328000
$ws.Cells.Item(4, 24).Value = "No es Simulación"
